# Correct the day/month ordering in the "Joining Date" (column L) text
# values for a few rows — the underlying strings were stored as plain
# text (not real dates), so we must keep writing them as text. Prefixing
# with a leading apostrophe forces Excel to store the literal text
# instead of auto-converting the date-shaped string into a date serial.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L3").Value = "'08/11/2020"
$ws.Range("L4").Value = "'08/11/2020"
$ws.Range("L6").Value = "'06/10/2020"
$ws.Range("L9").Value = "'04/05/2021"
